$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "97.662.71"
$ws.Range("E2").Value = "  +5.66%  "

Set-TextValue $ws.Range("D3") "3.126.68"
$ws.Range("E3").Value = "  +0.50%  "

$ws.Range("E4").Value = "  -0.04%  "

Set-TextValue $ws.Range("D5") "241.12"
$ws.Range("E5").Value = "  +2.84%  "

Set-TextValue $ws.Range("D6") "611.43"
$ws.Range("E6").Value = "  -0.24%  "

Set-TextValue $ws.Range("D7") "1.11"
$ws.Range("E7").Value = "  +3.43%  "

Set-TextValue $ws.Range("D8") "0.387"
$ws.Range("E8").Value = "  -0.40%  "

$ws.Range("E9").Value = "  +0.00%  "

Set-TextValue $ws.Range("D10") "3.122.52"
$ws.Range("E10").Value = "  +0.43%  "

Set-TextValue $ws.Range("D11") "0.790"
$ws.Range("E11").Value = "  -0.33%  "

$ws.Range("E12").Value = "  -0.12%  "

Set-TextValue $ws.Range("D13") "97.066.65"
$ws.Range("E13").Value = "  +5.28%  "

Set-TextValue $ws.Range("D14") "0.0000242"
$ws.Range("E14").Value = "  -0.73%  "

Set-TextValue $ws.Range("D15") "34.06"
$ws.Range("E15").Value = "  +0.66%  "

$ws.Range("E16").Value = "  -0.84%  "

Set-TextValue $ws.Range("D17") "3.710.03"
$ws.Range("E17").Value = "  +0.52%  "

Set-TextValue $ws.Range("D18") "3.125.87"
$ws.Range("E18").Value = "  +1.81%  "

$ws.Range("E19").Value = "  -5.21%  "

Set-TextValue $ws.Range("D20") "513.05"
$ws.Range("E20").Value = "  +16.94%  "

Set-TextValue $ws.Range("D21") "14.63"
$ws.Range("E21").Value = "  +0.37%  "

Set-TextValue $ws.Range("D22") "5.71"
$ws.Range("E22").Value = "  -1.89%  "

Set-TextValue $ws.Range("D23") "0.0000194"
$ws.Range("E23").Value = "  -5.30%  "

Set-TextValue $ws.Range("D24") "8.89"
$ws.Range("E24").Value = "  -3.88%  "

Set-TextValue $ws.Range("D25") "5.54"
$ws.Range("E25").Value = "  -0.65%  "

Set-TextValue $ws.Range("D26") "86.60"
$ws.Range("E26").Value = "  +1.47%  "

Set-TextValue $ws.Range("D27") "11.65"
$ws.Range("E27").Value = "  +1.45%  "

Set-TextValue $ws.Range("D28") "3.290.14"
$ws.Range("E28").Value = "  +0.67%  "

Set-TextValue $ws.Range("D29") "0.999"
$ws.Range("E29").Value = "  -0.02%  "

Set-TextValue $ws.Range("D30") "0.239"
$ws.Range("E30").Value = "  +4.78%  "

$ws.Range("E31").Value = "  -1.14%  "

Set-TextValue $ws.Range("D32") "0.126"
$ws.Range("E32").Value = "  +6.28%  "

Set-TextValue $ws.Range("D33") "9.09"
$ws.Range("E33").Value = "  -1.23%  "

Set-TextValue $ws.Range("D34") "26.49"
$ws.Range("E34").Value = "  +2.47%  "

Set-TextValue $ws.Range("D35") "0.840"
$ws.Range("E35").Value = "  -19.01%  "

Set-TextValue $ws.Range("D36") "0.153"
$ws.Range("E36").Value = "  -3.19%  "

Set-TextValue $ws.Range("D37") "7.38"
$ws.Range("E37").Value = "  -7.64%  "

Set-TextValue $ws.Range("D38") "489.36"
$ws.Range("E38").Value = "  +4.98%  "

Set-TextValue $ws.Range("D39") "1.88"
$ws.Range("E39").Value = "  -1.07%  "

Set-TextValue $ws.Range("D40") "24.25"
$ws.Range("E40").Value = "  +1.60%  "

Set-TextValue $ws.Range("D41") "0.439"
$ws.Range("E41").Value = "  +1.73%  "

$ws.Range("E42").Value = "  -2.40%  "

Set-TextValue $ws.Range("D43") "3.62"
$ws.Range("E43").Value = "  -7.83%  "

$ws.Range("E45").Value = "  -3.03%  "

Set-TextValue $ws.Range("D46") "163.16"
$ws.Range("E46").Value = "  +2.16%  "

Set-TextValue $ws.Range("D47") "1.94"
$ws.Range("E47").Value = "  +5.72%  "

Set-TextValue $ws.Range("D48") "0.697"
$ws.Range("E48").Value = "  +2.03%  "

Set-TextValue $ws.Range("D49") "44.48"
$ws.Range("E49").Value = "  +1.50%  "

# Row 50/51 swap: VeChain moves to row 50, Filecoin moves to row 51
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D50") "0.0327"
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D51") "4.40"
$ws.Range("E51").Value = "  +1.28%  "
